# Update the routes/stores coverage worksheet:
#  - add four new store columns (E:H) with their own route-coverage data
#  - add a new route row (row 16, "Z396") for the "Makro Comas" store
#  - rename "Plaza Vea Universitaria" (old D1) -> split into new stores
#  - split old route "Z423" into "Z423A" (row 14) and "Z423B" (row 11)
#  - update coverage flags across existing rows 2-15 for the new columns

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Header row (store names) ----
$ws.Range("B1").Value = "Hipermercado Metro Independencia"
$ws.Range("C1").Value = "Plaza Vea Izaguirre"
$ws.Range("D1").Value = "Tottus Mega Plaza"
$ws.Range("E1").Value = "Plaza Vea Los Olivos"
$ws.Range("F1").Value = "Tottus Los Olivos"
$ws.Range("G1").Value = "Makro Plaza Lima Norte"
$ws.Range("H1").Value = "Makro Comas"

# ---- Route labels (column A) ----
$labels = @("Z408","Z651","Z407","Z414","Z409","Z403","Z412","Z405","Z417","Z423B","Z411","Z399","Z423A","Z398","Z396")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}

# ---- Coverage matrix (rows 2-16, columns B-H) ----
$values = @(
    @(1, 1, 1, 0, 0, 1, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(1, 1, 1, 0, 0, 1, 0),
    @(1, 0, 1, 1, 0, 1, 0),
    @(1, 0, 1, 1, 0, 1, 0),
    @(0, 0, 0, 0, 1, 1, 0),
    @(1, 1, 1, 0, 0, 1, 0),
    @(1, 1, 1, 0, 0, 1, 0),
    @(0, 0, 1, 1, 0, 1, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 1, 1, 0),
    @(0, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 1, 1, 1),
    @(0, 0, 0, 0, 0, 0, 1)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    for ($c = 0; $c -lt 7; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $values[$r][$c]
    }
}

# ---- Formatting ----
# Propagate the existing header/label style (bold, centered, bordered) to the
# newly added header cells (E1:H1) and the newly added label cell (A16),
# without creating any new style definitions.
$ws.Range("B1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$excel.CutCopyMode = $false
